# Scheduled runner update: refresh market-price-derived columns
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)
# across the per-job profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 318.75
$ws.Range("I33").Value = 328.57144
$ws.Range("K33").Value = 328.57144
$ws.Range("M33").Value = -99.57144
$ws.Range("H55").Value = 453.26828
$ws.Range("I55").Value = 459.15384
$ws.Range("J55").Value = 443.06668
$ws.Range("K55").Value = 459.15384
$ws.Range("L55").Value = 443.06668
$ws.Range("M55").Value = -245.15384
$ws.Range("N55").Value = -871.06668
$ws.Range("H98").Value = 695067.3
$ws.Range("I98").Value = 794330.9399999999
$ws.Range("K98").Value = 794330.9399999999
$ws.Range("M98").Value = -792832.9399999999
$ws.Range("H122").Value = 695067.3
$ws.Range("I122").Value = 794330.9399999999
$ws.Range("K122").Value = 2382992.82
$ws.Range("M122").Value = -2380542.82
$ws.Range("H129").Value = 168889.08
$ws.Range("I129").Value = 252569
$ws.Range("J129").Value = 1529.25
$ws.Range("K129").Value = 757707
$ws.Range("L129").Value = 4587.75
$ws.Range("M129").Value = -752707
$ws.Range("N129").Value = -14587.75
$ws.Range("H132").Value = 1314.8387
$ws.Range("I132").Value = 1221.7307
$ws.Range("K132").Value = 3665.1921
$ws.Range("M132").Value = -1135.1921
$ws.Range("H138").Value = 2369.82
$ws.Range("J138").Value = 2494.4065
$ws.Range("L138").Value = 7483.2195
$ws.Range("N138").Value = -17763.2195
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13578.108
$ws.Range("I32").Value = 8660.709000000001
$ws.Range("K32").Value = 8660.709000000001
$ws.Range("M32").Value = -8373.709000000001
$ws.Range("H61").Value = 4316.231
$ws.Range("I61").Value = 2388.9375
$ws.Range("J61").Value = 7399.9
$ws.Range("K61").Value = 2388.9375
$ws.Range("L61").Value = 7399.9
$ws.Range("M61").Value = -2176.9375
$ws.Range("N61").Value = -7823.9
$ws.Range("H74").Value = 2129.4866
$ws.Range("I74").Value = 1839.258
$ws.Range("K74").Value = 1839.258
$ws.Range("M74").Value = -965.258
$ws.Range("H77").Value = 2129.4866
$ws.Range("I77").Value = 1839.258
$ws.Range("K77").Value = 9196.290000000001
$ws.Range("M77").Value = -4828.290000000001
$ws.Range("H97").Value = 844.3
$ws.Range("J97").Value = 835.4286
$ws.Range("L97").Value = 835.4286
$ws.Range("N97").Value = -1827.4286
$ws.Range("H122").Value = 3624.9
$ws.Range("I122").Value = 2749.92
$ws.Range("K122").Value = 8249.76
$ws.Range("M122").Value = -5799.76
$ws.Range("H132").Value = 4570.8335
$ws.Range("I132").Value = 4192.9414
$ws.Range("K132").Value = 12578.8242
$ws.Range("M132").Value = -10048.8242
$ws.Range("H136").Value = 4316.231
$ws.Range("I136").Value = 2388.9375
$ws.Range("J136").Value = 7399.9
$ws.Range("K136").Value = 7166.8125
$ws.Range("L136").Value = 22199.7
$ws.Range("M136").Value = -4616.8125
$ws.Range("N136").Value = -27299.7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 5856.4
$ws.Range("I82").Value = 5856.4
$ws.Range("K82").Value = 5856.4
$ws.Range("M82").Value = -5473.4
$ws.Range("H85").Value = 5856.4
$ws.Range("I85").Value = 5856.4
$ws.Range("K85").Value = 5856.4
$ws.Range("M85").Value = -4530.4
$ws.Range("H86").Value = 2105.5186
$ws.Range("J86").Value = 1547
$ws.Range("L86").Value = 1547
$ws.Range("N86").Value = -3793
$ws.Range("H89").Value = 2105.5186
$ws.Range("J89").Value = 1547
$ws.Range("L89").Value = 7735
$ws.Range("N89").Value = -18967
$ws.Range("H94").Value = 7694453
$ws.Range("I94").Value = 1654.8
$ws.Range("J94").Value = 33337114
$ws.Range("K94").Value = 1654.8
$ws.Range("L94").Value = 33337114
$ws.Range("M94").Value = -1203.8
$ws.Range("N94").Value = -33338016
$ws.Range("H97").Value = 12000
$ws.Range("I97").Value = 12000
$ws.Range("K97").Value = 12000
$ws.Range("M97").Value = -11009
$ws.Range("H105").Value = 6025
$ws.Range("I105").Value = 5037.875
$ws.Range("K105").Value = 5037.875
$ws.Range("M105").Value = -3290.875
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 50652.5
$ws.Range("I18").Value = 50621
$ws.Range("J18").Value = 50684
$ws.Range("K18").Value = 50621
$ws.Range("L18").Value = 50684
$ws.Range("M18").Value = -50391
$ws.Range("N18").Value = -51144
$ws.Range("H31").Value = 4688.185
$ws.Range("J31").Value = 5107.8237
$ws.Range("L31").Value = 5107.8237
$ws.Range("N31").Value = -5697.8237
$ws.Range("H34").Value = 4688.185
$ws.Range("J34").Value = 5107.8237
$ws.Range("L34").Value = 5107.8237
$ws.Range("N34").Value = -5511.8237
$ws.Range("H94").Value = 8308.799999999999
$ws.Range("I94").Value = 17904
$ws.Range("K94").Value = 17904
$ws.Range("M94").Value = -17453
$ws.Range("H132").Value = 3761.0476
$ws.Range("I132").Value = 3166
$ws.Range("K132").Value = 9498
$ws.Range("M132").Value = -6968
$ws.Range("H134").Value = 2747.2856
$ws.Range("I134").Value = 1415.8462
$ws.Range("K134").Value = 4247.5386
$ws.Range("M134").Value = -1712.5386
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1582.7587
$ws.Range("I129").Value = 1056.5
$ws.Range("J129").Value = 1720.0435
$ws.Range("K129").Value = 3169.5
$ws.Range("L129").Value = 5160.1305
$ws.Range("M129").Value = 1830.5
$ws.Range("N129").Value = -15160.1305
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 15929437
$ws.Range("I80").Value = 68578.17999999999
$ws.Range("K80").Value = 68578.17999999999
$ws.Range("M80").Value = -67580.17999999999
$ws.Range("H83").Value = 15929437
$ws.Range("I83").Value = 68578.17999999999
$ws.Range("K83").Value = 342890.9
$ws.Range("M83").Value = -337898.9
$ws.Range("H113").Value = 4632.353
$ws.Range("I113").Value = 2998
$ws.Range("K113").Value = 2998
$ws.Range("M113").Value = -828
$ws.Range("H132").Value = 4216.2354
$ws.Range("I132").Value = 3478.9048
$ws.Range("J132").Value = 7657.1113
$ws.Range("K132").Value = 10436.7144
$ws.Range("L132").Value = 22971.3339
$ws.Range("M132").Value = -7906.714399999999
$ws.Range("N132").Value = -28031.3339
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 4306291
$ws.Range("I52").Value = 6012808
$ws.Range("K52").Value = 6012808
$ws.Range("M52").Value = -6012582
$ws.Range("H62").Value = 37039988
$ws.Range("J62").Value = 55558056
$ws.Range("L62").Value = 55558056
$ws.Range("N62").Value = -55559304
$ws.Range("H65").Value = 37039988
$ws.Range("J65").Value = 55558056
$ws.Range("L65").Value = 277790280
$ws.Range("N65").Value = -277796520
$ws.Range("H122").Value = 3081.5
$ws.Range("I122").Value = 2522.8572
$ws.Range("K122").Value = 7568.571599999999
$ws.Range("M122").Value = -5118.571599999999
$ws.Range("H132").Value = 1760.2709
$ws.Range("I132").Value = 1060.8049
$ws.Range("K132").Value = 3182.4147
$ws.Range("M132").Value = -652.4147000000003
$ws.Range("H136").Value = 11367541
$ws.Range("I136").Value = 18521462
$ws.Range("K136").Value = 55564386
$ws.Range("M136").Value = -55561836
